$d = $word.ActiveDocument

# 1) "Firebase" -> "AJAX" in the technical capabilities list
$d.Content.Find.Execute("Firebase", $true, $false, $false, $false, $false,
                         $true, 1, $false, "AJAX", 2)

# 2) Move the "_GoBack" bookmark: Word keeps only a single "_GoBack"
#    bookmark, so adding a new one at the empty paragraph right after the
#    summary text automatically removes the old one that sat after
#    "History of producing well under pressure in a fast-paced sales
#    environment".
$summaryBlank = $d.Paragraphs.Item(8)
$d.Bookmarks.Add("_GoBack", $summaryBlank.Range)
